$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-survey record needs to be inserted before the current row 522,
# pushing row 522 (and everything below it) down by one row. The new row
# reuses the same market/category/unit metadata as the record that used to
# sit at row 522 ("$/bandeja 18 kilos" from "Región de Arica y Parinacota"),
# but carries its own date/volume/price figures.
$ws.Rows("522:522").Insert()

$ws.Range("A522").Value = 5
$ws.Range("B522").Value = "Macroferia Regional de Talca"
$ws.Range("C522").Value = "Maule"
$ws.Range("D522").Value = 44694
$ws.Range("E522").Value = 7
$ws.Range("F522").Value = 100112020
$ws.Range("G522").Value = "Tomate"
$ws.Range("H522").Value = "Larga vida"
$ws.Range("I522").Value = "Primera"
$ws.Range("J522").Value = 1500
$ws.Range("K522").Value = 22000
$ws.Range("L522").Value = 22000
$ws.Range("M522").Value = 22000
$ws.Range("N522").Value = "`$/bandeja 18 kilos"
$ws.Range("O522").Value = "Región de Arica y Parinacota"
$ws.Range("P522").Value = 1222
$ws.Range("Q522").Value = 18
$ws.Range("R522").Value = "Hortaliza"
